$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.9
$ws.Range("D3").Value = -7.382
$ws.Range("E5").Value = 13.098
$ws.Range("D14").Value = -8.077
$ws.Range("D21").Value = -7.9
$ws.Range("D23").Value = -7.505
$ws.Range("D25").Value = -8.307000000000002
